$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "new_name" column (F) with sequential renamed bmp filenames
$ws.Range("F1").Value = "new_name"
$ws.Range("F2").Value = "1.bmp"
$ws.Range("F3").Value = "2.bmp"
$ws.Range("F4").Value = "3.bmp"
$ws.Range("F5").Value = "4.bmp"
$ws.Range("F6").Value = "5.bmp"
$ws.Range("F7").Value = "6.bmp"
$ws.Range("F8").Value = "7.bmp"
$ws.Range("F9").Value = "8.bmp"
$ws.Range("F10").Value = "9.bmp"
$ws.Range("F11").Value = "10.bmp"
$ws.Range("F12").Value = "11.bmp"
$ws.Range("F13").Value = "12.bmp"
$ws.Range("F14").Value = "13.bmp"
$ws.Range("F15").Value = "14.bmp"
$ws.Range("F16").Value = "15.bmp"
$ws.Range("F17").Value = "16.bmp"

# Widen column E (gender) to fit content
$ws.Columns.Item(5).ColumnWidth = 23.5

# Swap out the U0361.bmp row for a new U0164.bmp entry (image moved out of stim folder)
$ws.Range("A13").Value = "U0164.bmp"
$ws.Range("B13").Value = 0.45
$ws.Range("D13").Value = 0.5

# Move selection to reflect where editing ended
$ws.Range("A14").Select() | Out-Null
